$d = $word.ActiveDocument

# Each entry: paragraph index (1-based, matching $d.Paragraphs), old text, new text.
# The Find is scoped to the specific paragraphs Range (Wrap = wdFindStop) so that
# duplicate text values elsewhere in the document do not interfere. The edit order is
# chosen (topologically) so that a paragraph whose *old* text equals another edits
# *new* text is processed first, avoiding any transient text collisions (e.g. "903÷2=").
$edits = @(
    ,@(1, "2024-04-05 Friday", "2024-04-06 Saturday")
    ,@(2, "538÷6=", "303÷8=")
    ,@(3, "989÷6=", "449÷3=")
    ,@(4, "710÷5=", "875÷8=")
    ,@(5, "974÷5=", "332÷9=")
    ,@(6, "757÷5=", "322÷6=")
    ,@(26, "601÷5=", "844÷8=")
    ,@(27, "524÷8=", "946÷3=")
    ,@(28, "730÷5=", "832÷4=")
    ,@(29, "942÷4=", "950÷5=")
    ,@(30, "650÷4=", "571÷3=")
    ,@(50, "393÷2=", "947÷9=")
    ,@(51, "545÷9=", "783÷8=")
    ,@(52, "416÷6=", "191÷6=")
    ,@(53, "221÷6=", "684÷2=")
    ,@(74, "869÷3=", "329÷2=")
    ,@(75, "825÷5=", "831÷4=")
    ,@(76, "903÷2=", "216÷7=")
    ,@(54, "718÷9=", "903÷2=")
    ,@(77, "250÷8=", "364÷4=")
    ,@(78, "851÷5=", "461÷9=")
    ,@(98, "189÷4=", "324÷7=")
    ,@(99, "655÷3=", "830÷2=")
    ,@(100, "420÷3=", "367÷6=")
    ,@(101, "850÷6=", "338÷6=")
    ,@(102, "629÷2=", "438÷8=")
)

foreach ($edit in $edits) {
    $paraIndex = $edit[0]
    $old = $edit[1]
    $new = $edit[2]
    $range = $d.Paragraphs.Item($paraIndex).Range
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 0, $false, $new, 2)
}

$d.Save()
